$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 178
$ws.Cells.Item(178, 1).Value = 237679041654
$ws.Cells.Item(178, 2).Value = "PIERRE MARIVOT TEMEZEU"
$ws.Cells.Item(178, 3).Value = "'0"
$ws.Cells.Item(178, 3).Style = "Normal"
$ws.Cells.Item(178, 4).Value = "Ndogbong 2"
$ws.Cells.Item(178, 5).Value = 10285
$ws.Cells.Item(178, 6).Value = 6876
$ws.Cells.Item(178, 7).Value = -3409
$ws.Cells.Item(178, 8).Value = 0.6685464268351969
$ws.Cells.Item(178, 9).Value = "Ndogbong"

# Row 179
$ws.Cells.Item(179, 1).Value = 237673671238
$ws.Cells.Item(179, 2).Value = "LA NEGRESSE SARL MBONE NDEMOU EPSE KAMSU ROSINE"
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = "Ndokoti Carrefour"
$ws.Cells.Item(179, 5).Value = 21090
$ws.Cells.Item(179, 6).Value = 4750
$ws.Cells.Item(179, 7).Value = -16340
$ws.Cells.Item(179, 8).Value = 0.2252252252252252
$ws.Cells.Item(179, 9).Value = "Ndogbong"

# Row 180
$ws.Cells.Item(180, 1).Value = 237652275301
$ws.Cells.Item(180, 2).Value = "NDAMI EPSE NONGA ROSALIE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(180, 4).Value = "Ndogbong Iut Ocm"
$ws.Cells.Item(180, 5).Value = 10000
$ws.Cells.Item(180, 6).Value = 7688
$ws.Cells.Item(180, 7).Value = -2312
$ws.Cells.Item(180, 8).Value = 0.7688
$ws.Cells.Item(180, 9).Value = "Ndogbong"

# Row 181
$ws.Cells.Item(181, 1).Value = 237681662701
$ws.Cells.Item(181, 2).Value = "TOUGOUA PAYOU JULIO OMER ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = "Bp Cite Ocm"
$ws.Cells.Item(181, 5).Value = 15300
$ws.Cells.Item(181, 6).Value = 17160
$ws.Cells.Item(181, 7).Value = 1860
$ws.Cells.Item(181, 8).Value = 1.12156862745098
$ws.Cells.Item(181, 9).Value = "Cite Sic"

# Row 182
$ws.Cells.Item(182, 1).Value = 237654349065
$ws.Cells.Item(182, 2).Value = "YASSI A BAA BELMOND CHIC MOBILE"
$ws.Cells.Item(182, 3).Value = "Rte_8"
$ws.Cells.Item(182, 4).Value = "Ndokoti Carrefour"
$ws.Cells.Item(182, 5).Value = 100420
$ws.Cells.Item(182, 6).Value = 360662
$ws.Cells.Item(182, 7).Value = 260242
$ws.Cells.Item(182, 8).Value = 3.591535550687114
$ws.Cells.Item(182, 9).Value = "Ndogbong"

# Row 183
$ws.Cells.Item(183, 1).Value = 237675637054
$ws.Cells.Item(183, 2).Value = "N A SOKOUDJOU DZOKOU"
$ws.Cells.Item(183, 3).Value = "Rte_5"
$ws.Cells.Item(183, 4).Value = "Ndogbong"
$ws.Cells.Item(183, 5).Value = 10000
$ws.Cells.Item(183, 6).Value = 6316
$ws.Cells.Item(183, 7).Value = -3684
$ws.Cells.Item(183, 8).Value = 0.6316
$ws.Cells.Item(183, 9).Value = "Ndogbong"

# Row 184
$ws.Cells.Item(184, 1).Value = 237671262234
$ws.Cells.Item(184, 2).Value = "DEUGOUE TOKO EPSE DIBANGUE LOISE LAURE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Cells.Item(184, 3).Value = "Rte_3"
$ws.Cells.Item(184, 4).Value = "Agape Ocm"
$ws.Cells.Item(184, 5).Value = 14940
$ws.Cells.Item(184, 6).Value = 70
$ws.Cells.Item(184, 7).Value = -14870
$ws.Cells.Item(184, 8).Value = 0.004685408299866131
$ws.Cells.Item(184, 9).Value = "Cite Sic"

# Row 185
$ws.Cells.Item(185, 1).Value = 237674890488
$ws.Cells.Item(185, 2).Value = "Sandrine Nkendji"
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = "Cite Bassa"
$ws.Cells.Item(185, 5).Value = 137280
$ws.Cells.Item(185, 6).Value = 253822
$ws.Cells.Item(185, 7).Value = 116542
$ws.Cells.Item(185, 8).Value = 1.84893648018648
$ws.Cells.Item(185, 9).Value = "Cite Sic"

# Row 186
$ws.Cells.Item(186, 1).Value = 237654079053
$ws.Cells.Item(186, 2).Value = "JUDITH AIMEE JOELEFACK JAZET EPSE NGUMATIO"
$ws.Cells.Item(186, 3).Value = 0
$ws.Cells.Item(186, 4).Value = "Ndokoti Carrefour"
$ws.Cells.Item(186, 5).Value = 13100
$ws.Cells.Item(186, 6).Value = 216137
$ws.Cells.Item(186, 7).Value = 203037
$ws.Cells.Item(186, 8).Value = 16.49900763358779
$ws.Cells.Item(186, 9).Value = "Ndogbong"

# Row 187
$ws.Cells.Item(187, 1).Value = 237652643069
$ws.Cells.Item(187, 2).Value = "CHIREL DELRICH TCHAPDA"
$ws.Cells.Item(187, 3).Value = 0
$ws.Cells.Item(187, 4).Value = "Pk8"
$ws.Cells.Item(187, 5).Value = 10000
$ws.Cells.Item(187, 6).Value = 6503
$ws.Cells.Item(187, 7).Value = -3497
$ws.Cells.Item(187, 8).Value = 0.6503
$ws.Cells.Item(187, 9).Value = "Ndogbong"

# Row 188
$ws.Cells.Item(188, 1).Value = 237673041651
$ws.Cells.Item(188, 2).Value = "DYLAN KEPSEU SIME"
$ws.Cells.Item(188, 3).Value = "Rte_1"
$ws.Cells.Item(188, 4).Value = "Ndokoti Carrefour"
$ws.Cells.Item(188, 5).Value = 39400
$ws.Cells.Item(188, 6).Value = 817810
$ws.Cells.Item(188, 7).Value = 778410
$ws.Cells.Item(188, 8).Value = 20.75659898477157
$ws.Cells.Item(188, 9).Value = "Ndogbong"

